$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I ("I0") and J ("IF"), reusing H1's style
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-25: column I is always 1, column J mirrors column H
for ($row = 2; $row -le 25; $row++) {
    $hValue = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 9).Value2 = 1
    $ws.Cells.Item($row, 10).Value2 = $hValue
}
